$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows at position 854, pushing existing rows 854-905 down
# to become rows 859-910.
$ws.Rows.Item(854).Insert()
$ws.Rows.Item(854).Insert()
$ws.Rows.Item(854).Insert()
$ws.Rows.Item(854).Insert()
$ws.Rows.Item(854).Insert()

# Populate the 5 new rows (854-858) with the new weekly price records.
$newRows = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44610, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103004, "Durazno", "Carson",       "Primera",  60, 13000, 13000, 13000, "`$/caja 15 kilos empedrada", "Región de O'Higgins", 867,  15),
    @(3, "Femacal de La Calera", "Coquimbo", 44610, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103004, "Durazno", "Carson",       "Segunda",  60, 11000, 11000, 11000, "`$/caja 15 kilos empedrada", "Región de O'Higgins", 733,  15),
    @(3, "Femacal de La Calera", "Coquimbo", 44610, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103004, "Durazno", "Doctor Davis", "Especial", 70, 16000, 16000, 16000, "`$/caja 15 kilos empedrada", "Región de O'Higgins", 1067, 15),
    @(3, "Femacal de La Calera", "Coquimbo", 44610, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103004, "Durazno", "Doctor Davis", "Primera",  78, 14000, 14000, 14000, "`$/caja 15 kilos empedrada", "Región de O'Higgins", 933,  15),
    @(3, "Femacal de La Calera", "Coquimbo", 44610, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103004, "Durazno", "Doctor Davis", "Segunda",  70, 12000, 12000, 12000, "`$/caja 15 kilos empedrada", "Región de O'Higgins", 800,  15)
)

$startRow = 854
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
